# Add 2022-Q3 data
#
# 1) "总计" (summary) sheet: insert a new row 2 for the 2022-Q3 quarter,
#    shifting the existing 2022-Q2 / 2022-Q1 / 2021-Q4 rows down by one
#    and renumbering the index column.
# 2) Duplicate the "2022-Q2" sheet (to inherit its layout/styling) to
#    create the new "2022-Q3" sheet positioned right before "2022-Q2",
#    then overwrite its holdings data with the 2022-Q3 figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Insert a fresh row above the current row 2 (old 2022-Q2 row).
$total.Rows.Item(2).Insert(-4121) | Out-Null

# The inserted row picked up stray formatting on B2:D2 - clear it so it
# matches the unstyled data cells used throughout the table.
$total.Range("B2:D2").ClearFormats() | Out-Null

# A2 lost the index-column style during the insert; restore it from A3
# (still carrying the original style) via a format-only paste.
$total.Range("A3").Copy() | Out-Null
$total.Range("A2").PasteSpecial(-4122) | Out-Null

# New 2022-Q3 row.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 1.86

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# ---------------------------------------------------------------------
# 2. New "2022-Q3" holdings sheet
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")

# Duplicate the 2022-Q2 sheet in place (placed immediately before it) so
# the new sheet starts with identical layout/styles, then rename it.
$q2.Copy($q2) | Out-Null
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Helper values: column A (index 0..9) already matches and is left as-is.
# Overwrite B..H for every data row with the 2022-Q3 figures. Numeric-
# looking text (fund codes / percentages) is written with a leading
# apostrophe so it stays literal text (preserves leading zeros / exact
# decimal formatting) instead of being coerced into a number.

# Row 2
$q3.Range("B2").Value = "'012526"
$q3.Range("C2").Value = "广发盛锦混合型证券投资基金A"
$q3.Range("D2").Value = "'25.18"
$q3.Range("E2").Value = "'90.19"
$q3.Range("F2").Value = "'3.99"
$q3.Range("G2").Value = "'1.0047"
$q3.Range("H2").Value = 7

# Row 3
$q3.Range("B3").Value = "'002446"
$q3.Range("C3").Value = "广发利鑫灵活配置混合A"
$q3.Range("D3").Value = "'13.77"
$q3.Range("E3").Value = "'74.30"
$q3.Range("F3").Value = "'4.37"
$q3.Range("G3").Value = "'0.6017"
$q3.Range("H3").Value = 2

# Row 4
$q3.Range("B4").Value = "'501030"
$q3.Range("C4").Value = "汇添富中证环境治理指数（LOF）A"
$q3.Range("D4").Value = "'3.12"
$q3.Range("E4").Value = "'92.74"
$q3.Range("F4").Value = "'2.08"
$q3.Range("G4").Value = "'0.0649"
$q3.Range("H4").Value = 9

# Row 5
$q3.Range("B5").Value = "'011172"
$q3.Range("C5").Value = "广发利鑫灵活配置混合C"
$q3.Range("D5").Value = "'1.41"
$q3.Range("E5").Value = "'74.30"
$q3.Range("F5").Value = "'4.37"
$q3.Range("G5").Value = "'0.0616"
$q3.Range("H5").Value = 2

# Row 6
$q3.Range("B6").Value = "'012527"
$q3.Range("C6").Value = "广发盛锦混合型证券投资基金C"
$q3.Range("D6").Value = "'1.16"
$q3.Range("E6").Value = "'90.19"
$q3.Range("F6").Value = "'3.99"
$q3.Range("G6").Value = "'0.0463"
$q3.Range("H6").Value = 7

# Row 7
$q3.Range("B7").Value = "'164908"
$q3.Range("C7").Value = "交银施罗德中证环境治理指数（LOF）"
$q3.Range("D7").Value = "'1.57"
$q3.Range("E7").Value = "'93.62"
$q3.Range("F7").Value = "'2.10"
$q3.Range("G7").Value = "'0.0330"
$q3.Range("H7").Value = 8

# Row 8
$q3.Range("B8").Value = "'501031"
$q3.Range("C8").Value = "汇添富中证环境治理指数（LOF）C"
$q3.Range("D8").Value = "'1.30"
$q3.Range("E8").Value = "'92.74"
$q3.Range("F8").Value = "'2.08"
$q3.Range("G8").Value = "'0.0270"
$q3.Range("H8").Value = 9

# Row 9
$q3.Range("B9").Value = "'006890"
$q3.Range("C9").Value = "上投摩根领先优选混合"
$q3.Range("D9").Value = "'0.32"
$q3.Range("E9").Value = "'80.46"
$q3.Range("F9").Value = "'3.15"
$q3.Range("G9").Value = "'0.0101"
$q3.Range("H9").Value = 7

# Row 10
$q3.Range("B10").Value = "'350007"
$q3.Range("C10").Value = "天治趋势精选混合"
$q3.Range("D10").Value = "'0.37"
$q3.Range("E10").Value = "'82.33"
$q3.Range("F10").Value = "'2.35"
$q3.Range("G10").Value = "'0.0087"
$q3.Range("H10").Value = 5

# Row 11
$q3.Range("B11").Value = "'013413"
$q3.Range("C11").Value = "交银施罗德中证环境治理指数（LOF）C"
$q3.Range("D11").Value = "'0.09"
$q3.Range("E11").Value = "'93.62"
$q3.Range("F11").Value = "'2.10"
$q3.Range("G11").Value = "'0.0019"
$q3.Range("H11").Value = 8
